$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 21: fill in issn, journal, publisher for the entry that previously
# only had id (A21) and color (C21).
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "Pediatric Dentistry"
$ws.Range("F21").Value = "Ingenta"

# Leave the selection where the author last left it when saving.
$ws.Range("E16").Select() | Out-Null
